$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "总计" (Total) sheet: insert a new row for the 2022-Q4 entry above the
#    existing 2022-Q3 row, shifting everything else down by one.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)

$total.Rows.Item(2).Insert()

# Column A is a plain running index (0,1,2,...); renumber the shifted rows.
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3

# The insert leaves the new row 2 without the bordered/centered style used by
# the rest of column A; copy it over from row 3 (which still carries the
# original formatting after the shift). Columns B:D stay unstyled, same as
# every other data row, so reset whatever Insert carried down from the
# header row.
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)
$total.Range("B2:D2").Style = "Normal"

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.58

# ---------------------------------------------------------------------------
# 2. Add the new "2022-Q4" detail sheet. Duplicating the existing "2022-Q3"
#    sheet keeps every style/number-format byte-identical, so only the cells
#    that actually changed need to be touched afterwards.
# ---------------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($q3)
$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

# Row 2 (fund 377016) - fund code/name stay the same, metrics are refreshed.
$q4.Range("D2").Value = "'25.54"
$q4.Range("D2").Style = "Normal"
$q4.Range("E2").Value = "'91.64"
$q4.Range("E2").Style = "Normal"
$q4.Range("F2").Value = "'2.15"
$q4.Range("F2").Style = "Normal"
$q4.Range("G2").Value = "'0.5491"
$q4.Range("G2").Style = "Normal"
$q4.Range("H2").Value = 7

# Row 3 (fund 006105) - fund code/name stay the same, metrics are refreshed.
$q4.Range("D3").Value = "'0.67"
$q4.Range("D3").Style = "Normal"
$q4.Range("E3").Value = "'91.81"
$q4.Range("E3").Style = "Normal"
$q4.Range("F3").Value = "'5.02"
$q4.Range("F3").Style = "Normal"
$q4.Range("G3").Value = "'0.0336"
$q4.Range("G3").Style = "Normal"
$q4.Range("H3").Value = 4
